# Update "想去人数" (want-to-go count) values in column F across the
# four sheets of the 广州-漫展信息 workbook, matching the re-scraped
# data published for gh-pages output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 22
$ws.Range("F3").Value  = 2747
$ws.Range("F4").Value  = 1081
$ws.Range("F5").Value  = 19937
$ws.Range("F7").Value  = 2315
$ws.Range("F8").Value  = 756
$ws.Range("F10").Value = 447
$ws.Range("F11").Value = 700
$ws.Range("F12").Value = 247
$ws.Range("F15").Value = 381
$ws.Range("F16").Value = 83
$ws.Range("F17").Value = 270
$ws.Range("F19").Value = 207

# Sheet 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value  = 20
$ws.Range("F7").Value  = 293
$ws.Range("F14").Value = 91
$ws.Range("F16").Value = 96

# Sheet 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6016
$ws.Range("F3").Value = 653
$ws.Range("F4").Value = 597

# Sheet 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 6016
$ws.Range("F3").Value  = 653
$ws.Range("F4").Value  = 597
$ws.Range("F6").Value  = 22
$ws.Range("F8").Value  = 2747
$ws.Range("F9").Value  = 1081
$ws.Range("F10").Value = 19937
$ws.Range("F12").Value = 20
$ws.Range("F15").Value = 293
$ws.Range("F16").Value = 2315
$ws.Range("F17").Value = 756
$ws.Range("F20").Value = 447
$ws.Range("F21").Value = 700
$ws.Range("F22").Value = 247
$ws.Range("F28").Value = 381
$ws.Range("F29").Value = 83
$ws.Range("F32").Value = 270
$ws.Range("F33").Value = 91
$ws.Range("F36").Value = 207
$ws.Range("F37").Value = 96
$ws.Range("F38").Value = 96
